$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15 (hunk 0)
$ws.Range("H15").Value = 266453.56
$ws.Range("I15").Value = 266453.56
$ws.Range("K15").Value = 799360.6799999999
$ws.Range("M15").Value = -799191.6799999999

# Row 43 (hunk 1)
$ws.Range("H43").Value = 871.125
$ws.Range("I43").Value = 855
$ws.Range("J43").Value = 898
$ws.Range("K43").Value = 855
$ws.Range("L43").Value = 898
$ws.Range("M43").Value = -786
$ws.Range("N43").Value = -1036

# Row 76 (hunk 2)
$ws.Range("H76").Value = 5167.0415
$ws.Range("I76").Value = 4476.6665
$ws.Range("K76").Value = 4476.6665
$ws.Range("M76").Value = -4161.6665

# Row 79 (hunk 3)
$ws.Range("H79").Value = 5167.0415
$ws.Range("I79").Value = 4476.6665
$ws.Range("K79").Value = 4476.6665
$ws.Range("M79").Value = -3384.6665

# Row 112 (hunk 4)
$ws.Range("H112").Value = 3370.5366
$ws.Range("J112").Value = 3507.7104
$ws.Range("L112").Value = 10523.1312
$ws.Range("N112").Value = -12739.1312

# Row 118 (hunk 5)
$ws.Range("H118").Value = 2930.4
$ws.Range("I118").Value = 1909
$ws.Range("J118").Value = 4462.5
$ws.Range("K118").Value = 5727
$ws.Range("L118").Value = 13387.5
$ws.Range("M118").Value = -4070
$ws.Range("N118").Value = -16701.5


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 6)
$ws.Range("H32").Value = 31263498
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 31263498
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 31263498
$ws.Range("N32").Value = -31264072
$ws.Range("M32").ClearContents()

# Row 45 (hunk 7)
$ws.Range("H45").Value = 4168.727
$ws.Range("I45").Value = 4032.4614
$ws.Range("K45").Value = 4032.4614
$ws.Range("M45").Value = -3655.4614

# Row 122 (hunk 8)
$ws.Range("H122").Value = 6118.1665
$ws.Range("I122").Value = 5541.8
$ws.Range("K122").Value = 16625.4
$ws.Range("M122").Value = -14175.4

# Row 128 (hunk 9)
$ws.Range("H128").Value = 99924.5
$ws.Range("J128").Value = 99924.5
$ws.Range("L128").Value = 99924.5
$ws.Range("N128").Value = -109884.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107 (hunk 10)
$ws.Range("H107").Value = 1985.3334
$ws.Range("I107").Value = 981.1429000000001
$ws.Range("K107").Value = 981.1429000000001
$ws.Range("M107").Value = 938.8570999999999


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 11)
$ws.Range("H16").Value = 3390.3
$ws.Range("I16").Value = 3238.5
$ws.Range("J16").Value = 3997.5
$ws.Range("K16").Value = 3238.5
$ws.Range("L16").Value = 3997.5
$ws.Range("M16").Value = -2951.5
$ws.Range("N16").Value = -4571.5

# Row 31 (hunk 12)
$ws.Range("H31").Value = 7205.135
$ws.Range("I31").Value = 3416.3103
$ws.Range("J31").Value = 11982.348
$ws.Range("K31").Value = 3416.3103
$ws.Range("L31").Value = 11982.348
$ws.Range("M31").Value = -3121.3103
$ws.Range("N31").Value = -12572.348

# Row 34 (hunk 13)
$ws.Range("H34").Value = 7205.135
$ws.Range("I34").Value = 3416.3103
$ws.Range("J34").Value = 11982.348
$ws.Range("K34").Value = 3416.3103
$ws.Range("L34").Value = 11982.348
$ws.Range("M34").Value = -3214.3103
$ws.Range("N34").Value = -12386.348

# Row 75 (hunk 14)
$ws.Range("H75").Value = 96497.164
$ws.Range("J75").Value = 102797
$ws.Range("L75").Value = 102797
$ws.Range("N75").Value = -104793

# Row 78 (hunk 15)
$ws.Range("H78").Value = 96497.164
$ws.Range("J78").Value = 102797
$ws.Range("L78").Value = 308391
$ws.Range("N78").Value = -318375

# Row 105 (hunk 16)
$ws.Range("H105").Value = 2755
$ws.Range("I105").Value = 2755
$ws.Range("K105").Value = 2755
$ws.Range("M105").Value = -1008

# Row 113 (hunk 17)
$ws.Range("H113").Value = 3390.3
$ws.Range("I113").Value = 3238.5
$ws.Range("J113").Value = 3997.5
$ws.Range("K113").Value = 3238.5
$ws.Range("L113").Value = 3997.5
$ws.Range("M113").Value = -1068.5
$ws.Range("N113").Value = -8337.5

# Row 124 (hunk 18)
$ws.Range("H124").Value = 23680.428
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()

# Row 132 (hunk 19)
$ws.Range("H132").Value = 1351.75
$ws.Range("I132").Value = 971.2632
$ws.Range("K132").Value = 2913.7896
$ws.Range("M132").Value = -383.7896000000001

# Row 133 (hunk 20)
$ws.Range("H133").Value = 20000
$ws.Range("I133").Value = 20000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 20000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -17470
$ws.Range("N133").ClearContents()

# Row 134 (hunk 21)
$ws.Range("H134").Value = 2213.6667
$ws.Range("I134").Value = 1881
$ws.Range("K134").Value = 5643
$ws.Range("M134").Value = -3108

# Row 135 (hunk 22)
$ws.Range("H135").Value = 78000
$ws.Range("J135").Value = 78000
$ws.Range("L135").Value = 78000
$ws.Range("N135").Value = -88140

# Row 141 (hunk 23)
$ws.Range("H141").Value = 610608.0600000001
$ws.Range("J141").Value = 654299.75
$ws.Range("L141").Value = 654299.75
$ws.Range("N141").Value = -664659.75


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122 (hunk 24)
$ws.Range("H122").Value = 2572.682
$ws.Range("I122").Value = 2148.2307
$ws.Range("J122").Value = 3185.7778
$ws.Range("K122").Value = 6444.6921
$ws.Range("L122").Value = 9557.3334
$ws.Range("M122").Value = -3994.6921
$ws.Range("N122").Value = -14457.3334

# Row 128 (hunk 25)
$ws.Range("H128").Value = 152740
$ws.Range("J128").Value = 152740
$ws.Range("L128").Value = 152740
$ws.Range("N128").Value = -162700


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40 (hunk 26)
$ws.Range("H40").Value = 4709.591
$ws.Range("I40").Value = 4780.55
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 4780.55
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -4644.55
$ws.Range("N40").Value = -4272

# Row 46 (hunk 27)
$ws.Range("H46").Value = 2231.8262
$ws.Range("I46").Value = 604.3077
$ws.Range("J46").Value = 4347.6
$ws.Range("K46").Value = 604.3077
$ws.Range("L46").Value = 4347.6
$ws.Range("M46").Value = -416.3077
$ws.Range("N46").Value = -4723.6

# Row 55 (hunk 28)
$ws.Range("H55").Value = 468.07407
$ws.Range("I55").Value = 390.375
$ws.Range("J55").Value = 581.0909
$ws.Range("K55").Value = 390.375
$ws.Range("L55").Value = 581.0909
$ws.Range("M55").Value = -217.375
$ws.Range("N55").Value = -927.0909

# Row 122 (hunk 29)
$ws.Range("H122").Value = 5483.9287
$ws.Range("J122").Value = 3865
$ws.Range("L122").Value = 11595
$ws.Range("N122").Value = -16495

# Row 128 (hunk 30)
$ws.Range("H128").Value = 84998.57000000001
$ws.Range("J128").Value = 84998.57000000001
$ws.Range("L128").Value = 84998.57000000001
$ws.Range("N128").Value = -94958.57000000001

# Row 136 (hunk 31)
$ws.Range("H136").Value = 4138.769
$ws.Range("I136").Value = 2989.5
$ws.Range("K136").Value = 8968.5
$ws.Range("M136").Value = -6418.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 92 (hunk 32)
$ws.Range("H92").Value = 44900
$ws.Range("J92").Value = 44900
$ws.Range("L92").Value = 44900
$ws.Range("N92").Value = -49892

# Row 122 (hunk 33)
$ws.Range("H122").Value = 1651.64
$ws.Range("I122").Value = 1686.6364
$ws.Range("K122").Value = 5059.9092
$ws.Range("M122").Value = -2609.9092

# Row 128 (hunk 34)
$ws.Range("H128").Value = 181792.25
$ws.Range("J128").Value = 181792.25
$ws.Range("L128").Value = 181792.25
$ws.Range("N128").Value = -191752.25

# Row 132 (hunk 35)
$ws.Range("H132").Value = 5730.8667
$ws.Range("I132").Value = 4799
$ws.Range("K132").Value = 14397
$ws.Range("M132").Value = -11867

